# Auto-generated edit script applying the diff to Ultros_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2169.95
$ws.Cells.Item(15, 9).Value = 2169.95
$ws.Cells.Item(15, 11).Value = 6509.849999999999
$ws.Cells.Item(15, 13).Value = -6340.849999999999
$ws.Cells.Item(19, 8).Value = 1032.7778
$ws.Cells.Item(19, 9).Value = 797.5
$ws.Cells.Item(19, 10).Value = 1100
$ws.Cells.Item(19, 11).Value = 797.5
$ws.Cells.Item(19, 12).Value = 1100
$ws.Cells.Item(19, 13).Value = -622.5
$ws.Cells.Item(19, 14).Value = -1450
$ws.Cells.Item(41, 9).Value = 200
$ws.Cells.Item(41, 11).Value = 200
$ws.Cells.Item(41, 13).Value = 240
$ws.Cells.Item(103, 8).Value = 699.8
$ws.Cells.Item(103, 9).Value = 649.5
$ws.Cells.Item(103, 11).Value = 1948.5
$ws.Cells.Item(103, 13).Value = -1362.5
$ws.Cells.Item(137, 8).Value = 5332.722
$ws.Cells.Item(137, 9).Value = 4900.5
$ws.Cells.Item(137, 10).Value = 5456.2144
$ws.Cells.Item(137, 11).Value = 14701.5
$ws.Cells.Item(137, 12).Value = 16368.6432
$ws.Cells.Item(137, 13).Value = -12151.5
$ws.Cells.Item(137, 14).Value = -21468.6432
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1913
$ws.Cells.Item(32, 9).Value = 2008.0625
$ws.Cells.Item(32, 10).Value = 1342.625
$ws.Cells.Item(32, 11).Value = 2008.0625
$ws.Cells.Item(32, 12).Value = 1342.625
$ws.Cells.Item(32, 13).Value = -1721.0625
$ws.Cells.Item(32, 14).Value = -1916.625
$ws.Cells.Item(45, 8).Value = 4732.875
$ws.Cells.Item(45, 9).Value = 3814.6667
$ws.Cells.Item(45, 11).Value = 3814.6667
$ws.Cells.Item(45, 13).Value = -3437.6667
$ws.Cells.Item(113, 8).Value = 199999
$ws.Cells.Item(113, 10).Value = 199999
$ws.Cells.Item(113, 12).Value = 199999
$ws.Cells.Item(113, 14).Value = -208677
$ws.Cells.Item(132, 8).Value = 2007
$ws.Cells.Item(132, 9).Value = 2007
$ws.Cells.Item(132, 11).Value = 6021
$ws.Cells.Item(132, 13).Value = -3491
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(70, 8).Value = 79999.164
$ws.Cells.Item(70, 10).Value = 79999.164
$ws.Cells.Item(70, 12).Value = 79999.164
$ws.Cells.Item(70, 14).Value = -80585.164
$ws.Cells.Item(73, 8).Value = 79999.164
$ws.Cells.Item(73, 10).Value = 79999.164
$ws.Cells.Item(73, 12).Value = 79999.164
$ws.Cells.Item(73, 14).Value = -82027.164
$ws.Cells.Item(134, 8).Value = 3670.8333
$ws.Cells.Item(134, 10).Value = 2338
$ws.Cells.Item(134, 12).Value = 7014
$ws.Cells.Item(134, 14).Value = -12084
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 1835681.6
$ws.Cells.Item(6, 9).Value = 2202318
$ws.Cells.Item(6, 11).Value = 2202318
$ws.Cells.Item(6, 13).Value = -2202205
$ws.Cells.Item(31, 8).Value = 2942.087
$ws.Cells.Item(31, 9).Value = 1799.6923
$ws.Cells.Item(31, 11).Value = 1799.6923
$ws.Cells.Item(31, 13).Value = -1504.6923
$ws.Cells.Item(34, 8).Value = 2942.087
$ws.Cells.Item(34, 9).Value = 1799.6923
$ws.Cells.Item(34, 11).Value = 1799.6923
$ws.Cells.Item(34, 13).Value = -1597.6923
$ws.Cells.Item(107, 8).Value = 5002.909
$ws.Cells.Item(107, 9).Value = 206.75
$ws.Cells.Item(107, 11).Value = 206.75
$ws.Cells.Item(107, 13).Value = 1713.25
$ws.Cells.Item(132, 8).Value = 3194
$ws.Cells.Item(132, 9).Value = 2622
$ws.Cells.Item(132, 11).Value = 7866
$ws.Cells.Item(132, 13).Value = -5336
$ws.Cells.Item(134, 8).Value = 7997.636
$ws.Cells.Item(134, 9).Value = 7997.636
$ws.Cells.Item(134, 11).Value = 23992.908
$ws.Cells.Item(134, 13).Value = -21457.908
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 306.4
$ws.Cells.Item(14, 9).Value = 306.4
$ws.Cells.Item(14, 11).Value = 919.1999999999999
$ws.Cells.Item(14, 13).Value = -746.1999999999999
$ws.Cells.Item(22, 8).Value = 799
$ws.Cells.Item(22, 9).Value = 799
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 2397
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -2228
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 799
$ws.Cells.Item(27, 9).Value = 799
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 2397
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -2295
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(44, 8).Value = 172.25
$ws.Cells.Item(44, 9).Value = 172.25
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 516.75
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = -118.75
$ws.Cells.Item(44, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 1179104
$ws.Cells.Item(131, 10).Value = 3988.889
$ws.Cells.Item(131, 12).Value = 11966.667
$ws.Cells.Item(131, 14).Value = -22046.667
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 5025500
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 5025500
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 5025500
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).Value = -5025724
$ws.Cells.Item(8, 8).Value = 5025500
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 5025500
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 5025500
$ws.Cells.Item(8, 13).ClearContents()
$ws.Cells.Item(8, 14).Value = -5025778
$ws.Cells.Item(43, 8).Value = 26565.13
$ws.Cells.Item(43, 9).Value = 15599.8
$ws.Cells.Item(43, 11).Value = 15599.8
$ws.Cells.Item(43, 13).Value = -15448.8
$ws.Cells.Item(70, 8).Value = 130594.78
$ws.Cells.Item(70, 9).Value = 229901.2
$ws.Cells.Item(70, 11).Value = 229901.2
$ws.Cells.Item(70, 13).Value = -229631.2
$ws.Cells.Item(73, 8).Value = 130594.78
$ws.Cells.Item(73, 9).Value = 229901.2
$ws.Cells.Item(73, 11).Value = 229901.2
$ws.Cells.Item(73, 13).Value = -228965.2
$ws.Cells.Item(80, 8).Value = 118190.6
$ws.Cells.Item(80, 9).Value = 374371.66
$ws.Cells.Item(80, 10).Value = 8398.714
$ws.Cells.Item(80, 11).Value = 374371.66
$ws.Cells.Item(80, 12).Value = 8398.714
$ws.Cells.Item(80, 13).Value = -373373.66
$ws.Cells.Item(80, 14).Value = -10394.714
$ws.Cells.Item(83, 8).Value = 118190.6
$ws.Cells.Item(83, 9).Value = 374371.66
$ws.Cells.Item(83, 10).Value = 8398.714
$ws.Cells.Item(83, 11).Value = 1871858.3
$ws.Cells.Item(83, 12).Value = 41993.57
$ws.Cells.Item(83, 13).Value = -1866866.3
$ws.Cells.Item(83, 14).Value = -51977.57
$ws.Cells.Item(86, 8).Value = 29999
$ws.Cells.Item(86, 10).Value = 29999
$ws.Cells.Item(86, 12).Value = 29999
$ws.Cells.Item(86, 14).Value = -32371
$ws.Cells.Item(89, 8).Value = 29999
$ws.Cells.Item(89, 10).Value = 29999
$ws.Cells.Item(89, 12).Value = 89997
$ws.Cells.Item(89, 14).Value = -101853
$ws.Cells.Item(126, 8).Value = 3756
$ws.Cells.Item(126, 9).Value = 3756
$ws.Cells.Item(126, 11).Value = 11268
$ws.Cells.Item(126, 13).Value = -8798
$ws.Cells.Item(132, 8).Value = 9804
$ws.Cells.Item(132, 9).Value = 9058.333000000001
$ws.Cells.Item(132, 11).Value = 27174.999
$ws.Cells.Item(132, 13).Value = -24644.999
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 1000
$ws.Cells.Item(5, 11).Value = 1000
$ws.Cells.Item(5, 13).Value = -887
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4901.5
$ws.Cells.Item(7, 9).Value = 2804
$ws.Cells.Item(7, 11).Value = 2804
$ws.Cells.Item(7, 13).Value = -2692
$ws.Cells.Item(58, 8).Value = 9955.5
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 4901.5
$ws.Cells.Item(126, 9).Value = 2804
$ws.Cells.Item(126, 11).Value = 8412
$ws.Cells.Item(126, 13).Value = -5942
$ws.Cells.Item(132, 8).Value = 2922.3845
$ws.Cells.Item(132, 9).Value = 2299.7
$ws.Cells.Item(132, 11).Value = 6899.099999999999
$ws.Cells.Item(132, 13).Value = -4369.099999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 44003
$ws.Cells.Item(24, 9).Value = 44003
$ws.Cells.Item(24, 11).Value = 44003
$ws.Cells.Item(24, 13).Value = -43773
